$wb = $excel.ActiveWorkbook

$timestamp = "2025-11-30 03:04:59"

$ws = $wb.Worksheets.Item(2)
for ($r = 2; $r -le 26; $r++) {
    $ws.Range("AA" + $r).Value = $timestamp
}

$ws.Range("C6").Value = 30
$ws.Range("D6").Value = 495
$ws.Range("E6").Value = 221
$ws.Range("F6").Value = 274
$ws.Range("G6").Value = 16.5
$ws.Range("H6").Value = 7.37
$ws.Range("I6").Value = 9.130000000000001
$ws.Range("J6").Value = 103
$ws.Range("K6").Value = 117
$ws.Range("W6").Value = 16

$ws.Range("C7").Value = 19
$ws.Range("D7").Value = 257
$ws.Range("E7").Value = 111
$ws.Range("F7").Value = 146
$ws.Range("G7").Value = 13.53
$ws.Range("H7").Value = 5.84
$ws.Range("I7").Value = 7.68
$ws.Range("J7").Value = 53
$ws.Range("K7").Value = 53
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 4
$ws.Range("V7").Value = 14

$ws.Range("C15").Value = 20
$ws.Range("D15").Value = 358
$ws.Range("E15").Value = 171
$ws.Range("F15").Value = 187
$ws.Range("G15").Value = 17.9
$ws.Range("H15").Value = 8.550000000000001
$ws.Range("I15").Value = 9.35
$ws.Range("J15").Value = 63
$ws.Range("K15").Value = 81

$ws.Range("C17").Value = 19
$ws.Range("D17").Value = 296
$ws.Range("E17").Value = 109
$ws.Range("F17").Value = 187
$ws.Range("G17").Value = 15.58
$ws.Range("H17").Value = 5.74
$ws.Range("I17").Value = 9.84
$ws.Range("J17").Value = 52
$ws.Range("K17").Value = 76
$ws.Range("W17").Value = 10

$ws.Range("C20").Value = 29
$ws.Range("D20").Value = 505
$ws.Range("E20").Value = 220
$ws.Range("F20").Value = 285
$ws.Range("G20").Value = 17.41
$ws.Range("H20").Value = 7.59
$ws.Range("I20").Value = 9.83
$ws.Range("J20").Value = 100
$ws.Range("K20").Value = 105
$ws.Range("L20").Value = 4
$ws.Range("M20").Value = 9
$ws.Range("V20").Value = 16

$ws.Range("C23").Value = 18
$ws.Range("D23").Value = 228
$ws.Range("E23").Value = 87
$ws.Range("F23").Value = 141
$ws.Range("G23").Value = 12.67
$ws.Range("H23").Value = 4.83
$ws.Range("I23").Value = 7.83
$ws.Range("J23").Value = 41
$ws.Range("K23").Value = 58

$ws = $wb.Worksheets.Item(3)
for ($r = 2; $r -le 26; $r++) {
    $ws.Range("AA" + $r).Value = $timestamp
}

$ws.Range("C8").Value = 28
$ws.Range("D8").Value = 431
$ws.Range("E8").Value = 164
$ws.Range("F8").Value = 267
$ws.Range("G8").Value = 15.39
$ws.Range("H8").Value = 5.86
$ws.Range("I8").Value = 9.539999999999999
$ws.Range("J8").Value = 77
$ws.Range("K8").Value = 106

$ws.Range("C11").Value = 20
$ws.Range("D11").Value = 280
$ws.Range("E11").Value = 124
$ws.Range("F11").Value = 156
$ws.Range("G11").Value = 14
$ws.Range("H11").Value = 6.2
$ws.Range("I11").Value = 7.8
$ws.Range("J11").Value = 62
$ws.Range("K11").Value = 73

$ws.Range("C14").Value = 27
$ws.Range("D14").Value = 446
$ws.Range("E14").Value = 225
$ws.Range("F14").Value = 221
$ws.Range("G14").Value = 16.52
$ws.Range("H14").Value = 8.33
$ws.Range("I14").Value = 8.19
$ws.Range("J14").Value = 110
$ws.Range("K14").Value = 103
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 3
$ws.Range("V14").Value = 12

$ws.Range("C18").Value = 31
$ws.Range("D18").Value = 517
$ws.Range("E18").Value = 239
$ws.Range("F18").Value = 278
$ws.Range("G18").Value = 16.68
$ws.Range("H18").Value = 7.71
$ws.Range("I18").Value = 8.970000000000001
$ws.Range("J18").Value = 112
$ws.Range("K18").Value = 119
$ws.Range("L18").Value = 3
$ws.Range("M18").Value = 4
$ws.Range("V18").Value = 8

$ws.Range("C24").Value = 31
$ws.Range("D24").Value = 546
$ws.Range("E24").Value = 217
$ws.Range("F24").Value = 329
$ws.Range("G24").Value = 17.61
$ws.Range("H24").Value = 7
$ws.Range("I24").Value = 10.61
$ws.Range("J24").Value = 96
$ws.Range("K24").Value = 127
$ws.Range("W24").Value = 20
